$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-formatted cells to remain text (they mirror price strings
# using "." as thousands separators, so Excel must not reinterpret them as numbers).
$textCells = @("D4","D5","D6","D9","D10","D11","D12","D13","D15","D16","D17","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D39","D40","D41","D42","D43","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.057.86"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "3.253.39"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "573.74"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").Value = "180.33"
$ws.Range("E6").Value = "  -4.12%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.248.42"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  -7.43%  "
$ws.Range("D11").Value = "0.561"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").Value = "45.52"
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "3.781.04"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").Value = "611.07"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "8.27"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "65.185.88"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "3.263.47"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "17.41"
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("D21").Value = "10.72"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "0.874"
$ws.Range("E22").Value = "  -3.90%  "
$ws.Range("D23").Value = "18.13"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").Value = "97.80"
$ws.Range("E25").Value = "  -5.13%  "
$ws.Range("D26").Value = "3.91"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "9.24"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").Value = "30.06"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "8.22"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("D31").Value = "6.35"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "544.59"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "3.54"
$ws.Range("E33").Value = "  -11.81%  "
$ws.Range("D34").Value = "10.67"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").Value = "3.738.62"
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("D36").Value = "0.101"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "55.47"
$ws.Range("E38").Value = "  -3.85%  "
$ws.Range("D39").Value = "0.125"
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("B40").Value = "ApeXProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D40").Value = "3.38"
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "32.00"
$ws.Range("E41").Value = "  -6.12%  "
$ws.Range("D42").Value = "3.08"
$ws.Range("E42").Value = "  -7.11%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.53"
$ws.Range("E43").Value = "  -6.76%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0663"
$ws.Range("E44").Value = "  -9.76%  "
$ws.Range("D45").Value = "0.324"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("D46").Value = "0.0399"
$ws.Range("E46").Value = "  -5.32%  "
$ws.Range("D47").Value = "3.01"
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "0.124"
$ws.Range("E49").Value = "  -3.92%  "
$ws.Range("D50").Value = "2.46"
$ws.Range("E50").Value = "  -5.67%  "
$ws.Range("D51").Value = "127.65"
$ws.Range("E51").Value = "  +4.36%  "
